# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted at row 23 (pushing the existing
# rows 23-50 down to 24-51), adding the latest "Camote" observation for
# Vega Central Mapocho de Santiago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23; this shifts rows 23:50 down to 24:51
# and carries their formatting (including the date style on column D).
$ws.Rows("23:23").Insert()

# Populate the newly inserted row 23 with the new weekly record.
$ws.Cells.Item(23, 1).Value  = 9
$ws.Cells.Item(23, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(23, 3).Value  = "Metropolitana"
$ws.Cells.Item(23, 4).Value  = 44494
$ws.Cells.Item(23, 5).Value  = 13
$ws.Cells.Item(23, 6).Value  = 100114002
$ws.Cells.Item(23, 7).Value  = "Camote"
$ws.Cells.Item(23, 8).Value  = "Sin especificar"
$ws.Cells.Item(23, 9).Value  = "Primera"
$ws.Cells.Item(23, 10).Value = 430
$ws.Cells.Item(23, 11).Value = 17000
$ws.Cells.Item(23, 12).Value = 17000
$ws.Cells.Item(23, 13).Value = 17000
$ws.Cells.Item(23, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(23, 15).Value = "Perú"
$ws.Cells.Item(23, 16).Value = 944
$ws.Cells.Item(23, 17).Value = 18
$ws.Cells.Item(23, 18).Value = "Hortaliza"
